$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B22 should become a genuine numeric value (currently stored as an inline string "4")
$ws.Range("B22").Value = 4

# New row 23 data
$ws.Range("A23").Value = "Ruilin"

# politeness_score "3" must stay text (not a number) for this row, like the rest of the row's
# columns - force text via the leading apostrophe, then clear the resulting quote-prefix style
# so the cell keeps the sheet's default formatting.
$ws.Range("B23").Value = "'3"
$ws.Range("B23").Style = "Normal"

$ws.Range("C23").Value = "无"
$ws.Range("D23").Value = "CRT"
$ws.Range("E23").Value = "WRI"
$ws.Range("F23").Value = "01473e7f-4f45-41be-bd3f-03c0ff83190e"
$ws.Range("G23").Value = "H1u8fMW0b_annotated.xlsx"
$ws.Range("H23").Value = "The citations are in non-standard format (section 1.2: Kalman (1960))."
